# Atualização de bases das ligas, do dia: 29-03-2024 às 17:05
#
# The source feed re-ordered a handful of fixtures, so the rows that used to
# hold two different matches now have their data swapped between the two
# rows (only the running index in column A stays put). A couple of other
# rows just received refreshed odds for a couple of columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    # Columns B..AC hold all the match data; column A is just the running
    # index and must stay where it is.
    $range1 = $ws.Range("B$row1" + ":AC$row1")
    $range2 = $ws.Range("B$row2" + ":AC$row2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

# Rows 73/74 (ids 7646749 / 7646750) swapped places.
Swap-RowData 73 74

# Rows 104/105 (ids 7127370 / 7127374) swapped places.
Swap-RowData 104 105

# Rows 124/125 (ids 7127388 / 7128012) swapped places.
Swap-RowData 124 125

# Row 133 (id 7126793): refreshed closing odds.
$ws.Range("R133").Value2 = 1.87
$ws.Range("S133").Value2 = 2.03

# Row 134 (id 7127396): refreshed odds.
$ws.Range("N134").Value2 = 2.2
$ws.Range("P134").Value2 = 3
$ws.Range("R134").Value2 = 1.93
$ws.Range("S134").Value2 = 1.97
